# DMS: Translate LuckyNumber Template
# Rename the "LuckyNumber" worksheet to the Vietnamese "Giải thưởng"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LuckyNumber")
$ws.Name = "Giải thưởng"
